$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293 - this pushes the existing rows
# 293..372 down to 294..373 (all column values, including dates/volumes,
# shift down by one row).
$ws.Rows.Item(293).Insert()

# Populate the newly inserted row 293 with the new weekly price record.
$ws.Range("A293").Value = 3
$ws.Range("B293").Value = "Femacal de La Calera"
$ws.Range("C293").Value = "Coquimbo"
$ws.Range("D293").Value = 44841
$ws.Range("E293").Value = 5
$ws.Range("F293").Value = 100112039
$ws.Range("G293").Value = "Ciboulette"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 120
$ws.Range("K293").Value = 1500
$ws.Range("L293").Value = 1500
$ws.Range("M293").Value = 1500
$ws.Range("N293").Value = "`$/docena de atados"
$ws.Range("O293").Value = "Provincia de Quillota"
$ws.Range("P293").Value = 500
$ws.Range("Q293").Value = 3
$ws.Range("R293").Value = "Hortaliza"
